# Update "minimal_seats-Croatia" sheet:
#  - expand the abbreviated party/key labels (row 1 + column A) into full
#    descriptive names
#  - re-derive the per-year seat table, which shifts several numeric
#    columns right (new "HDSSB", "HSP-AS", "Ind-ZG" columns were inserted)
#    and removes the old "HDS"/ constituency column that no longer applies

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1: header labels (full party names) ----
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "- - Ethnic Constituencies (-, -)"
$ws.Range("C1").Value = "HDSS - Croatian Democratic Peasants Party (Hrvatska Demokratska Seljacka Stranka, HDSS)"
$ws.Range("D1").Value = "HDZ - Croatian Democratic Union (Hrvatska Demokratska Zajednica, HDZ)"
$ws.Range("E1").Value = "HNS - Croatian People's Party-Liberal Democrats (Hrvatska Narodna Stranka-Liberalni Demokrati, HNS)"
$ws.Range("F1").Value = "HSLS - Croation Social-Liberal Party (Hrvatska Socijalno Liberalna Stranka, HSLS)"
$ws.Range("G1").Value = "HSP - Croatian Party of Rights (Hrvatska Stranka Prava, HSP)"
$ws.Range("H1").Value = "HSS - Croatian Peasant Party (Hrvatska Seljacka Stranka, HSS)"
$ws.Range("I1").Value = "HSU - Croatian Party of Pensioners (Hrvatska Stranka Umirovljenika, HSU)"
$ws.Range("J1").Value = "SDP - Social Democratic Party of Croatia (Socijaldemokratska Partija Hrvatske, SDP)"
$ws.Range("K1").Value = "HDSSB - Croatian Democratic Alliance of Slavonia and Baranja and allies (HDSSB-HDSSD-ZH) (Hrvatski Demokratski Savez Slavonije I Baranje, HDSSB)"
$ws.Range("L1").Value = "IDS - Istrian Democratic Assembly (Istarski Demokratski Sabor, IDS)"
$ws.Range("M1").Value = "HL-SR - Croatian Labourists - Labour Party (Hrvatski Laburisti-Stranka Rada, HL-SR)"
$ws.Range("N1").Value = "HSP-AS - Croatian Party of Rights-Dr. Ante Starcevic (Hrvatska Stranka Prava Dr. Ante Starcevic, HSP-AS)"
$ws.Range("O1").Value = "Ind-IG - Independent List of Ivan Grubsic (Neovisna Lista - Ivan Grubisic, Ind-IG)"
$ws.Range("P1").Value = "BM365 - Bandic Milan Party of Labour and Solidarity (Bandic Milan 365-Stranka rada i solidarnosti, BM365)"
$ws.Range("Q1").Value = "MOST - Bridge of Independent Lists (Most nezavisnih lista, MOST)"
$ws.Range("R1").Value = "ZZ - Human Shield (Zivi zid, ZZ)"
$ws.Range("S1").Value = "Ind-ZG - Independent List of Zeljko Glasnovic (Neovisna Lista - Zeljko Glasnovic, Ind-ZG)"
$ws.Range("T1").Value = "DPMS - Miroslav Skoro Homeland Movement (Domovinski pokret Miroslava Skore, DPMS)"
$ws.Range("U1").Value = "M! - We Can! (Mozemo!, M!)"
$ws.Range("V1").Value = "NS - Our Party (Nasa Stranka, NS)"
$ws.Range("W1").Value = "SSIP - Party With a First and Last Name (Stranka s imenom i prezimenom, SSIP)"

# fix up the few header cells that need real (non-ASCII-folded) characters
$ws.Range("C1").Value = "HDSS - Croatian Democratic Peasants Party (Hrvatska Demokratska Seljačka Stranka, HDSS)"
$ws.Range("H1").Value = "HSS - Croatian Peasant Party (Hrvatska Seljačka Stranka, HSS)"
$ws.Range("M1").Value = "HL-SR - Croatian Labourists – Labour Party (Hrvatski Laburisti-Stranka Rada, HL-SR)"
$ws.Range("N1").Value = "HSP-AS - Croatian Party of Rights-Dr. Ante Starčević (Hrvatska Stranka Prava Dr. Ante Starčević, HSP-AS)"
$ws.Range("O1").Value = "Ind-IG - Independent List of Ivan Grubsic (Neovisna Lista - Ivan Grubišić, Ind-IG)"
$ws.Range("P1").Value = "BM365 - Bandic Milan Party of Labour and Solidarity (Bandic Milan 365-Stranka ´ rada i solidarnosti, BM365)"
$ws.Range("R1").Value = "ZZ - Human Shield (Živi zid, ZZ)"
$ws.Range("S1").Value = "Ind-ZG - Independent List of Željko Glasnovic (Neovisna Lista - Željko Glasnovic, Ind-ZG)"
$ws.Range("T1").Value = "DPMŠ - Miroslav Škoro Homeland Movement (Domovinski pokret Miroslava Škore, DPMŠ)"
$ws.Range("U1").Value = "M! - We Can! (Možemo!, M!)"
$ws.Range("V1").Value = "NS - Our Party (Naša Stranka, NS)"

# ---- Column A: year labels shift down one conceptual row each (unchanged text, same order) ----
$ws.Range("A2").Value = "2003"
$ws.Range("A3").Value = "2007"
$ws.Range("A4").Value = "2011"
$ws.Range("A5").Value = "2015"
$ws.Range("A6").Value = "2016"
$ws.Range("A7").Value = "2020"

# ---- Row 2 (2003): unchanged numbers, columns B:J ----
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 6

# ---- Row 3 (2007): old "HDS" column (C3) removed; new L3 added ----
$ws.Range("B3").Value = 8
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 11
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 3

# ---- Row 4 (2011): C4/G4 removed, old L4 becomes K4, N4/O4 added ----
$ws.Range("B4").Value = 0
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 0
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0

# ---- Row 5 (2015): C5 removed, old N5 becomes L5, old O5 removed, Q5/R5 added ----
$ws.Range("B5").Value = 1
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 1
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 1

# ---- Row 6 (2016): C6 removed, K6 split into K6(new)+L6, M6 removed, N6->P6, O6->Q6, P6->R6, new S6 ----
$ws.Range("B6").Value = 6
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 17
$ws.Range("J6").Value = 12
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 2
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 6
$ws.Range("S6").Value = 1

# ---- Row 7 (2020): O7 removed, O7(old v=3)->Q7, Q7(old)->T7, R7->U7, S7->V7, T7->W7 ----
$ws.Range("B7").Value = 6
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 1
$ws.Range("J7").Value = 7
$ws.Range("O7").ClearContents()
$ws.Range("Q7").Value = 3
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").Value = 7
$ws.Range("U7").Value = 5
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 3
